$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = "Lonay – Préverenges – Vullierens"
$ws.Range("C10").Select()
